$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "consequents_length"
# Copy the header formatting from G1 (bold, centered, bordered) onto H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "consequents_length"

# Data values H2:H6
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 2
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
